$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "jo"
$ws.Range("A2").Value = "ahmed"
$ws.Range("B2").Select()
